# edit.ps1 - applies the "cleaned Main and MainWindow" changes to System explanation.docx
$d = $word.ActiveDocument

# --- 1. Paragraph 2 (originally empty) gets the "Main class" explanation + the _GoBack bookmark ---
$p2 = $d.Paragraphs.Item(2)
$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>The main class is, as it name indicates, Main. Please execute it to see the GUI that we have designed.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$p2.Range.InsertXML($xml2)

# --- 2. Paragraph "Administrator View" loses its lastRenderedPageBreak marker ---
$pAdminView = $d.Paragraphs.Item(16)
$xml16 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Administrator View</w:t></w:r></w:p>
'@
$pAdminView.Range.InsertXML($xml16)

# --- 3. "soyadmin" is split into two runs ("soyadmi" + "n") inside the Password paragraph ---
$pPassword = $d.Paragraphs.Item(18)
$xml18 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Password: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>soyadmi</w:t></w:r><w:r><w:t>n</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$pPassword.Range.InsertXML($xml18)

# --- 4. Remove the empty paragraph that used to follow the admin password paragraph ---
$pEmpty = $d.Paragraphs.Item(19)
$pEmpty.Range.Delete()

# --- 5. Replace the 4 admin-view description paragraphs with their rewritten / expanded versions ---
$pStart = $d.Paragraphs.Item(20)
$pEnd = $d.Paragraphs.Item(23)
$rBlock = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$xmlBlock = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The admin view is only </w:t></w:r><w:r><w:t>accessible</w:t></w:r><w:r><w:t xml:space="preserve"> if you are an employee with the administrator status, not just a regular employee.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The admin view´s main target is to facilitate a user-friendly graphical interface for the supervisors of the warehouses in order they can keep record of everything that goes in and out of their facilities </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>and also</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> see all of the drugs, employees, warehouses, providers and clients that the company owns. The supervisor can even see the picture of a specific drug or employee just be clicking on it.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The supervisor will not only be able </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>to  see</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> all of this data but to modify (by double clicking in some cells, not all </w:t></w:r><w:r><w:t>POJOs</w:t></w:r><w:r><w:t xml:space="preserve"> are modifiable), add or delete some of this data.</w:t></w:r><w:r><w:t xml:space="preserve"> When viewing employees or drugs, if you click on one of them, their picture, if there is one, will show up on the right. In either case if you click this one, you will be able to change the picture of the drug or employee.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>Finally</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the Admin View also counts with the option to generate both xml or html files from the database.</w:t></w:r><w:r><w:t xml:space="preserve"> The XML works perfectly, while the HTML has still a bug, that we don’t fully understand, because it doesn’t let us apply the XSLT to generate the database unless we have saved the XML by hand. Still, we have the HTML done, and it works, if you want to check it out with a version of the database that we have saved by hand (with this we mean that we have created it with Java, opened it, made a change, undone the change, and then save it) and you can see it works. If you execute the class </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>XMLManager</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> it will generate the HTML. It should also open automatically in your computer, but if not, the error will give you the direction to check the HTML.</w:t></w:r></w:p>
'@
$rBlock.InsertXML($xmlBlock)

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
